# Weekly fruit/vegetable price update: a new daily record is inserted
# before the existing row 485 (Perejil, Mercado Mayorista Lo Valledor de
# Santiago), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 485; Excel shifts rows 485:508 down to 486:509
# and copies formatting (incl. the date style on column D) from the row
# above, matching the target dimension A1:R509.
$ws.Rows("485:485").Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A485").Value = 6
$ws.Range("B485").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C485").Value = "Metropolitana"
$ws.Range("D485").Value = 44706
$ws.Range("E485").Value = 13
$ws.Range("F485").Value = 100112044
$ws.Range("G485").Value = "Perejil"
$ws.Range("H485").Value = "Sin especificar"
$ws.Range("I485").Value = "Primera"
$ws.Range("J485").Value = 210
$ws.Range("K485").Value = 8000
$ws.Range("L485").Value = 9000
$ws.Range("M485").Value = 8429
$ws.Range("N485").Value = "`$/docena de atados"
$ws.Range("O485").Value = "Región Metropolitana"
$ws.Range("P485").Value = 2810
$ws.Range("Q485").Value = 3
$ws.Range("R485").Value = "Hortaliza"
